# Applies the "created unified requirements doc" edit:
#  1. Rename the sheet "Sheet1" -> "Sprint 1"
#  2. Update the "Meeting 2" / Arpit "what will they do" text in C36
#  3. Move the saved view/selection: scroll so row 28 is at the top, and
#     select H38 (mirrors the saved sheetView's topLeftCell + selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sprint 1"

$ws.Range("C36").Value = "What will they do: Learn how to do use cases, install UMLet, install Android Studio, install ZenHub"

$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("H38").Select()
